$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.08032186512941
$ws.Range("C2").Value = 4.707088164900111
$ws.Range("E2").Value = 11.26418622786683
$ws.Range("F2").Value = 47.14247703679082
$ws.Range("G2").Value = 3.744802456999608
$ws.Range("I2").Value = 33.88940456219572
$ws.Range("J2").Value = 9.954836322397528
$ws.Range("K2").Value = 15.70608465718566
$ws.Range("L2").Value = 11.89756668590153
$ws.Range("N2").Value = 23.43051152018387
$ws.Range("B3").Value = 17.93516884377031
$ws.Range("C3").Value = 4.55538783772402
$ws.Range("E3").Value = 11.2693233596374
$ws.Range("F3").Value = 47.1093776534086
$ws.Range("G3").Value = 3.747594357931122
$ws.Range("I3").Value = 33.94169114173044
$ws.Range("J3").Value = 9.969337384346815
$ws.Range("K3").Value = 15.61164198458322
$ws.Range("L3").Value = 11.89444654952093
$ws.Range("N3").Value = 23.4887549136079
$ws.Range("B4").Value = 17.84988857762897
$ws.Range("C4").Value = 4.461241061494007
$ws.Range("E4").Value = 11.27386134642715
$ws.Range("F4").Value = 47.09867620329401
$ws.Range("G4").Value = 3.749398697020055
$ws.Range("I4").Value = 33.97935404845763
$ws.Range("J4").Value = 9.978629083806206
$ws.Range("K4").Value = 15.55707617801378
$ws.Range("L4").Value = 11.89441014164147
$ws.Range("N4").Value = 23.52645583803323
$ws.Range("B5").Value = 17.8161372069183
$ws.Range("C5").Value = 4.422691873472539
$ws.Range("E5").Value = 11.27605955736496
$ws.Range("F5").Value = 47.0967383937513
$ws.Range("G5").Value = 3.750156713615337
$ws.Range("I5").Value = 33.99609785204241
$ws.Range("J5").Value = 9.982513446563143
$ws.Range("K5").Value = 15.53571996263454
$ws.Range("L5").Value = 11.89486928576165
$ws.Range("N5").Value = 23.54230767187784
$ws.Range("B6").Value = 17.81059424433713
$ws.Range("C6").Value = 4.416281805820096
$ws.Range("E6").Value = 11.27644567146014
$ws.Range("F6").Value = 47.09656305005133
$ws.Range("G6").Value = 3.750283957027867
$ws.Range("I6").Value = 33.99896239399104
$ws.Range("J6").Value = 9.983164366623075
$ws.Range("K6").Value = 15.5322274502718
$ws.Range("L6").Value = 11.8949741906536
$ws.Range("N6").Value = 23.54496937765151
$ws.Range("B7").Value = 17.84942929773632
$ws.Range("C7").Value = 4.460721823461134
$ws.Range("E7").Value = 11.27388957814028
$ws.Range("F7").Value = 47.09864025455467
$ws.Range("G7").Value = 3.749408827746469
$ws.Range("I7").Value = 33.9795742122263
$ws.Range("J7").Value = 9.97868107274256
$ws.Range("K7").Value = 15.55678457407014
$ws.Range("L7").Value = 11.89441441298232
$ws.Range("N7").Value = 23.52666764290874
$ws.Range("B8").Value = 18.02949639445668
$ws.Range("C8").Value = 4.655029513673917
$ws.Range("E8").Value = 11.26567091934201
$ws.Range("F8").Value = 47.12907001084476
$ws.Range("G8").Value = 3.745746450624315
$ws.Range("I8").Value = 33.90627816035636
$ws.Range("J8").Value = 9.959755977528728
$ws.Range("K8").Value = 15.67282245660197
$ws.Range("L8").Value = 11.89610169261604
$ws.Range("N8").Value = 23.45019181103023
$ws.Range("B9").Value = 18.41151577947544
$ws.Range("C9").Value = 5.025347674606017
$ws.Range("E9").Value = 11.26048869337326
$ws.Range("F9").Value = 47.26486725504513
$ws.Range("G9").Value = 3.739275924723443
$ws.Range("I9").Value = 33.80673107603165
$ws.Range("J9").Value = 9.925705846798108
$ws.Range("K9").Value = 15.92661929441236
$ws.Range("L9").Value = 11.9142542870765
$ws.Range("N9").Value = 23.31557469636355
$ws.Range("B10").Value = 18.70756918467477
$ws.Range("C10").Value = 5.287531832254872
$ws.Range("E10").Value = 11.26328419725585
$ws.Range("F10").Value = 47.41066528976499
$ws.Range("G10").Value = 3.734950776724386
$ws.Range("I10").Value = 33.76063165769536
$ws.Range("J10").Value = 9.902532174420783
$ws.Range("K10").Value = 16.12778387662159
$ws.Range("L10").Value = 11.93653901061562
$ws.Range("N10").Value = 23.22598027287936
$ws.Range("B11").Value = 18.84509145547098
$ws.Range("C11").Value = 5.404051502486851
$ws.Range("E11").Value = 11.2659748164478
$ws.Range("F11").Value = 47.4868744691413
$ws.Range("G11").Value = 3.733075202861867
$ws.Range("I11").Value = 33.74554768871757
$ws.Range("J11").Value = 9.892384966339625
$ws.Range("K11").Value = 16.22219152416781
$ws.Range("L11").Value = 11.94859429184375
$ws.Range("N11").Value = 23.1872323465703
$ws.Range("B12").Value = 18.89753368692933
$ws.Range("C12").Value = 5.447735259319284
$ws.Range("E12").Value = 11.26719631706322
$ws.Range("F12").Value = 47.51714209235722
$ws.Range("G12").Value = 3.732378115233296
$ws.Range("I12").Value = 33.74068317690902
$ws.Range("J12").Value = 9.888598846179571
$ws.Range("K12").Value = 16.25833131838744
$ws.Range("L12").Value = 11.9534325628424
$ws.Range("N12").Value = 23.17284770008585
$ws.Range("B13").Value = 18.88622375733181
$ws.Range("C13").Value = 5.438347393268539
$ws.Range("E13").Value = 11.2669242568435
$ws.Range("F13").Value = 47.51056097541106
$ws.Range("G13").Value = 3.732527661750904
$ws.Range("I13").Value = 33.74169313366737
$ws.Range("J13").Value = 9.889411751523241
$ws.Range("K13").Value = 16.25053107192362
$ws.Range("L13").Value = 11.9523784474268
$ws.Range("N13").Value = 23.17593287495378
$ws.Range("B14").Value = 18.84939883549428
$ws.Range("C14").Value = 5.407654441276998
$ws.Range("E14").Value = 11.26607125760911
$ws.Range("F14").Value = 47.48933644660693
$ws.Range("G14").Value = 3.733017589878017
$ws.Range("I14").Value = 33.74513049369992
$ws.Range("J14").Value = 9.892072351396781
$ws.Range("K14").Value = 16.2251571137456
$ws.Range("L14").Value = 11.94898687956116
$ws.Range("N14").Value = 23.18604313707732
$ws.Range("B15").Value = 18.82688880584019
$ws.Range("C15").Value = 5.388795622409759
$ws.Range("E15").Value = 11.26557511479194
$ws.Range("F15").Value = 47.47651888535918
$ws.Range("G15").Value = 3.733319395588709
$ws.Range("I15").Value = 33.74734636316845
$ws.Range("J15").Value = 9.893709381549559
$ws.Range("K15").Value = 16.20966476793573
$ws.Range("L15").Value = 11.94694494447192
$ws.Range("N15").Value = 23.19227350076102
$ws.Range("B16").Value = 18.69863532319223
$ws.Range("C16").Value = 5.279857773447586
$ws.Range("E16").Value = 11.2631367919849
$ws.Range("F16").Value = 47.40588271080023
$ws.Range("G16").Value = 3.73507519415224
$ws.Range("I16").Value = 33.76173597584373
$ws.Range("J16").Value = 9.903203230058846
$ws.Range("K16").Value = 16.121670123808
$ws.Range("L16").Value = 11.93578954519415
$ws.Range("N16").Value = 23.22855292081625
$ws.Range("B17").Value = 18.62065395221688
$ws.Range("C17").Value = 5.212290686517872
$ws.Range("E17").Value = 11.26200339757793
$ws.Range("F17").Value = 47.3650730163393
$ws.Range("G17").Value = 3.736175820944964
$ws.Range("I17").Value = 33.772072062051
$ws.Range("J17").Value = 9.909128227675767
$ws.Range("K17").Value = 16.06841130668947
$ws.Range("L17").Value = 11.92943550053161
$ws.Range("N17").Value = 23.25132333665987
$ws.Range("B18").Value = 18.57607167170143
$ws.Range("C18").Value = 5.173171456508756
$ws.Range("E18").Value = 11.26148514548958
$ws.Range("F18").Value = 47.34253162718634
$ws.Range("G18").Value = 3.736817532121048
$ws.Range("I18").Value = 33.77857116864448
$ws.Range("J18").Value = 9.91257329236737
$ws.Range("K18").Value = 16.03805257007118
$ws.Range("L18").Value = 11.92596150689381
$ws.Range("N18").Value = 23.26460940769579
$ws.Range("B19").Value = 18.56102466836818
$ws.Range("C19").Value = 5.159883746444624
$ws.Range("E19").Value = 11.26133267005978
$ws.Range("F19").Value = 47.33505980411766
$ws.Range("G19").Value = 3.737036294020357
$ws.Range("I19").Value = 33.78086678246327
$ws.Range("J19").Value = 9.913746125131521
$ws.Range("K19").Value = 16.0278215473165
$ws.Range("L19").Value = 11.92481638059741
$ws.Range("N19").Value = 23.26914034327853
$ws.Range("B20").Value = 18.62892751764949
$ws.Range("C20").Value = 5.219510194322511
$ws.Range("E20").Value = 11.2621102274881
$ws.Range("F20").Value = 47.36932098138389
$ws.Range("G20").Value = 3.736057761612983
$ws.Range("I20").Value = 33.77091441686616
$ws.Range("J20").Value = 9.908493657796594
$ws.Range("K20").Value = 16.07405259957415
$ws.Range("L20").Value = 11.93009321732982
$ws.Range("N20").Value = 23.24887981672898
$ws.Range("B21").Value = 18.86020563771646
$ws.Range("C21").Value = 5.416681969109445
$ws.Range("E21").Value = 11.26631631677293
$ws.Range("F21").Value = 47.49553247704697
$ws.Range("G21").Value = 3.732873329754111
$ws.Range("I21").Value = 33.7440978518667
$ws.Range("J21").Value = 9.891289340522997
$ws.Range("K21").Value = 16.23259970258768
$ws.Range("L21").Value = 11.94997567237307
$ws.Range("N21").Value = 23.18306568644087
$ws.Range("B22").Value = 19.01346935224644
$ws.Range("C22").Value = 5.542964848977498
$ws.Range("E22").Value = 11.27024555852757
$ws.Range("F22").Value = 47.58622398357453
$ws.Range("G22").Value = 3.730868743352589
$ws.Range("I22").Value = 33.73151155530729
$ws.Range("J22").Value = 9.880373954049372
$ws.Range("K22").Value = 16.33847775764992
$ws.Range("L22").Value = 11.96456106832845
$ws.Range("N22").Value = 23.14173284822337
$ws.Range("B23").Value = 18.93149108653257
$ws.Range("C23").Value = 5.475814840121185
$ws.Range("E23").Value = 11.26804093813692
$ws.Range("F23").Value = 47.53707401628708
$ws.Range("G23").Value = 3.731931641386046
$ws.Range("I23").Value = 33.7377768638801
$ws.Range("J23").Value = 9.886169742757053
$ws.Range("K23").Value = 16.28177106017169
$ws.Range("L23").Value = 11.95663189259565
$ws.Range("N23").Value = 23.16363937007373
$ws.Range("B24").Value = 18.62518625585928
$ws.Range("C24").Value = 5.216247104918238
$ws.Range("E24").Value = 11.26206151422855
$ws.Range("F24").Value = 47.36739760759284
$ws.Range("G24").Value = 3.736111108363756
$ws.Range("I24").Value = 33.77143605397126
$ws.Range("J24").Value = 9.908780426234552
$ws.Range("K24").Value = 16.07150135766018
$ws.Range("L24").Value = 11.92979530571473
$ws.Range("N24").Value = 23.24998392439203
$ws.Range("B25").Value = 18.30530076029356
$ws.Range("C25").Value = 4.92668034673481
$ws.Range("E25").Value = 11.26072621837575
$ws.Range("F25").Value = 47.22001385208929
$ws.Range("G25").Value = 3.74095072639525
$ws.Range("I25").Value = 33.82891937758533
$ws.Range("J25").Value = 9.934591961457905
$ws.Range("K25").Value = 15.85528039895518
$ws.Range("L25").Value = 11.90776403364023
$ws.Range("N25").Value = 23.35035393215418
